# Adds season-record columns (Wins / Losses / Ties) to the roster table.
# Season for this team: 83 wins, 79 losses, 0 ties - same record repeated
# for every player row, matching the source data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, one column past the existing "Unnamed: 28" (AC) column.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold / centered / bordered header style used by the rest of row 1
# (copy formatting only from A1, so the header text we just set is preserved).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the team's season record for every player row (2-50).
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 83  # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 79  # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF -> Ties
}
